$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 971.7
$ws.Range("I19").Value = 358.33334
$ws.Range("K19").Value = 358.33334
$ws.Range("M19").Value = -183.33334
$ws.Range("H112").Value = 1002.3333
$ws.Range("J112").Value = 1007.8421
$ws.Range("L112").Value = 3023.5263
$ws.Range("N112").Value = -5239.5263
$ws.Range("H129").Value = 2483.9365
$ws.Range("J129").Value = 963.14813
$ws.Range("L129").Value = 2889.44439
$ws.Range("N129").Value = -12889.44439
$ws.Range("H132").Value = 8936262
$ws.Range("I132").Value = 9623557
$ws.Range("J132").Value = 1428
$ws.Range("K132").Value = 28870671
$ws.Range("L132").Value = 4284
$ws.Range("M132").Value = -28868141
$ws.Range("N132").Value = -9344
$ws.Range("H138").Value = 4614.94
$ws.Range("I138").Value = 2388.5
$ws.Range("J138").Value = 5867.3125
$ws.Range("K138").Value = 7165.5
$ws.Range("L138").Value = 17601.9375
$ws.Range("M138").Value = -2025.5
$ws.Range("N138").Value = -27881.9375

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26034.314
$ws.Range("I32").Value = 4286
$ws.Range("K32").Value = 4286
$ws.Range("M32").Value = -3999
$ws.Range("H45").Value = 1742.4
$ws.Range("I45").Value = 1632.4667
$ws.Range("J45").Value = 1907.3
$ws.Range("K45").Value = 1632.4667
$ws.Range("L45").Value = 1907.3
$ws.Range("M45").Value = -1255.4667
$ws.Range("N45").Value = -2661.3
$ws.Range("H63").Value = 2742.2222
$ws.Range("I63").Value = 2100
$ws.Range("J63").Value = 3545
$ws.Range("K63").Value = 2100
$ws.Range("L63").Value = 3545
$ws.Range("M63").Value = -1414
$ws.Range("N63").Value = -4917
$ws.Range("H66").Value = 2742.2222
$ws.Range("I66").Value = 2100
$ws.Range("J66").Value = 3545
$ws.Range("K66").Value = 10500
$ws.Range("L66").Value = 17725
$ws.Range("M66").Value = -7068
$ws.Range("N66").Value = -24589
$ws.Range("H119").Value = 34028
$ws.Range("J119").Value = 34028
$ws.Range("L119").Value = 34028
$ws.Range("N119").Value = -43704
$ws.Range("H122").Value = 1810.0333
$ws.Range("I122").Value = 1696.375
$ws.Range("K122").Value = 5089.125
$ws.Range("M122").Value = -2639.125

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 15809.875
$ws.Range("I82").Value = 3723.2
$ws.Range("J82").Value = 35954.332
$ws.Range("K82").Value = 3723.2
$ws.Range("L82").Value = 35954.332
$ws.Range("M82").Value = -3340.2
$ws.Range("N82").Value = -36720.332
$ws.Range("H85").Value = 15809.875
$ws.Range("I85").Value = 3723.2
$ws.Range("J85").Value = 35954.332
$ws.Range("K85").Value = 3723.2
$ws.Range("L85").Value = 35954.332
$ws.Range("M85").Value = -2397.2
$ws.Range("N85").Value = -38606.332
$ws.Range("H86").Value = 81166.86
$ws.Range("J86").Value = 3122.8333
$ws.Range("L86").Value = 3122.8333
$ws.Range("N86").Value = -5368.8333
$ws.Range("H89").Value = 81166.86
$ws.Range("J89").Value = 3122.8333
$ws.Range("L89").Value = 15614.1665
$ws.Range("N89").Value = -26846.1665
$ws.Range("H99").Value = 1203.05
$ws.Range("I99").Value = 1121.25
$ws.Range("J99").Value = 1257.5834
$ws.Range("K99").Value = 1121.25
$ws.Range("L99").Value = 1257.5834
$ws.Range("M99").Value = 376.75
$ws.Range("N99").Value = -4253.5834

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16422.25
$ws.Range("I31").Value = 33205.676
$ws.Range("J31").Value = 2360.4595
$ws.Range("K31").Value = 33205.676
$ws.Range("L31").Value = 2360.4595
$ws.Range("M31").Value = -32910.676
$ws.Range("N31").Value = -2950.4595
$ws.Range("H34").Value = 16422.25
$ws.Range("I34").Value = 33205.676
$ws.Range("J34").Value = 2360.4595
$ws.Range("K34").Value = 33205.676
$ws.Range("L34").Value = 2360.4595
$ws.Range("M34").Value = -33003.676
$ws.Range("N34").Value = -2764.4595
$ws.Range("H107").Value = 4199.2856
$ws.Range("I107").Value = 6866.0625
$ws.Range("J107").Value = 643.5833
$ws.Range("K107").Value = 6866.0625
$ws.Range("L107").Value = 643.5833
$ws.Range("M107").Value = -4946.0625
$ws.Range("N107").Value = -4483.5833

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 248953.48
$ws.Range("I107").Value = 385.7143
$ws.Range("K107").Value = 1157.1429
$ws.Range("M107").Value = 762.8571000000002
$ws.Range("H122").Value = 7457.4
$ws.Range("I122").Value = 266
$ws.Range("J122").Value = 10072.454
$ws.Range("K122").Value = 2394
$ws.Range("L122").Value = 90652.086
$ws.Range("M122").Value = 56
$ws.Range("N122").Value = -95552.086
$ws.Range("H131").Value = 1005.3917
$ws.Range("I131").Value = 748
$ws.Range("J131").Value = 1019.38043
$ws.Range("K131").Value = 2244
$ws.Range("L131").Value = 3058.14129
$ws.Range("M131").Value = 2796
$ws.Range("N131").Value = -13138.14129

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5099.875
$ws.Range("I43").Value = 1933.3334
$ws.Range("J43").Value = 6999.8
$ws.Range("K43").Value = 1933.3334
$ws.Range("L43").Value = 6999.8
$ws.Range("M43").Value = -1782.3334
$ws.Range("N43").Value = -7301.8
$ws.Range("H46").Value = 10450
$ws.Range("I46").Value = 9900
$ws.Range("K46").Value = 9900
$ws.Range("M46").Value = -9744
$ws.Range("H57").Value = 16500
$ws.Range("J57").Value = 16500
$ws.Range("L57").Value = 16500
$ws.Range("N57").Value = -18140
$ws.Range("H80").Value = 3317.8
$ws.Range("I80").Value = 5300
$ws.Range("J80").Value = 1996.3334
$ws.Range("K80").Value = 5300
$ws.Range("L80").Value = 1996.3334
$ws.Range("M80").Value = -4302
$ws.Range("N80").Value = -3992.3334
$ws.Range("H83").Value = 3317.8
$ws.Range("I83").Value = 5300
$ws.Range("J83").Value = 1996.3334
$ws.Range("K83").Value = 26500
$ws.Range("L83").Value = 9981.666999999999
$ws.Range("M83").Value = -21508
$ws.Range("N83").Value = -19965.667

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7001859.5
$ws.Range("I16").Value = 12600929
$ws.Range("J16").Value = 3022.75
$ws.Range("K16").Value = 12600929
$ws.Range("L16").Value = 3022.75
$ws.Range("M16").Value = -12600759
$ws.Range("N16").Value = -3362.75
$ws.Range("H36").Value = 41066.668
$ws.Range("J36").Value = 41066.668
$ws.Range("L36").Value = 41066.668
$ws.Range("N36").Value = -42190.668
$ws.Range("H119").Value = 39793.332
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 39793.332
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 39793.332
$ws.Range("M119").ClearContents()
$ws.Range("N119").Value = -49469.332
$ws.Range("H132").Value = 7792.8
$ws.Range("I132").Value = 8683.615
$ws.Range("J132").Value = 2002.5
$ws.Range("K132").Value = 26050.845
$ws.Range("L132").Value = 6007.5
$ws.Range("M132").Value = -23520.845
$ws.Range("N132").Value = -11067.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 7340.6665
$ws.Range("J74").Value = 8335.444
$ws.Range("L74").Value = 8335.444
$ws.Range("N74").Value = -10207.444
$ws.Range("H77").Value = 7340.6665
$ws.Range("J77").Value = 8335.444
$ws.Range("L77").Value = 25006.332
$ws.Range("N77").Value = -34366.33199999999
$ws.Range("H119").Value = 22398.666
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()
$ws.Range("H132").Value = 3766.697
$ws.Range("I132").Value = 3790.0715
$ws.Range("J132").Value = 3635.8
$ws.Range("K132").Value = 11370.2145
$ws.Range("L132").Value = 10907.4
$ws.Range("M132").Value = -8840.2145
$ws.Range("N132").Value = -15967.4

